$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'26.087.88"
$ws.Cells.Item(2, 4).ClearFormats()
$ws.Cells.Item(2, 5).Value = '  -0.98%  '
$ws.Cells.Item(3, 4).Value = "'1.666.05"
$ws.Cells.Item(3, 4).ClearFormats()
$ws.Cells.Item(3, 5).Value = '  -1.23%  '
$ws.Cells.Item(4, 4).Value = "'1.004"
$ws.Cells.Item(4, 4).ClearFormats()
$ws.Cells.Item(4, 5).Value = '  -0.72%  '
$ws.Cells.Item(5, 4).Value = "'209.63"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = '  -4.04%  '
$ws.Cells.Item(6, 4).Value = "'0.5162"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = '  -5.44%  '
$ws.Cells.Item(7, 5).Value = '  -0.63%  '
$ws.Cells.Item(8, 4).Value = "'0.2633"
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).Value = '  -3.27%  '
$ws.Cells.Item(9, 4).Value = "'0.06195"
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).Value = '  -3.84%  '
$ws.Cells.Item(10, 4).Value = "'20.90"
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).Value = '  -5.06%  '
$ws.Cells.Item(11, 4).Value = "'0.07488"
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).Value = '  -2.25%  '
$ws.Cells.Item(12, 4).Value = "'1.668.84"
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).Value = '  -1.16%  '
$ws.Cells.Item(13, 4).Value = "'4.422"
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).Value = '  -2.42%  '
$ws.Cells.Item(14, 4).Value = "'0.5570"
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).Value = '  -4.06%  '
$ws.Cells.Item(15, 4).Value = "'0.000007924"
$ws.Cells.Item(15, 4).ClearFormats()
$ws.Cells.Item(15, 5).Value = '  -4.91%  '
$ws.Cells.Item(16, 4).Value = "'65.31"
$ws.Cells.Item(16, 4).ClearFormats()
$ws.Cells.Item(16, 5).Value = '  +0.34%  '
$ws.Cells.Item(17, 4).Value = "'26.104.48"
$ws.Cells.Item(17, 4).ClearFormats()
$ws.Cells.Item(17, 5).Value = '  -1.07%  '
$ws.Cells.Item(18, 5).Value = '  -0.64%  '
$ws.Cells.Item(19, 4).Value = "'4.782"
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).Value = '  -3.24%  '
$ws.Cells.Item(20, 4).Value = "'10.38"
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).Value = '  -5.30%  '
$ws.Cells.Item(21, 4).Value = "'185.04"
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(21, 5).Value = '  -2.78%  '
$ws.Cells.Item(22, 4).Value = "'6.146"
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).Value = '  -1.22%  '
$ws.Cells.Item(23, 4).Value = "'1.005"
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(23, 5).Value = '  -0.67%  '
$ws.Cells.Item(24, 4).Value = "'145.99"
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).Value = '  -2.58%  '
$ws.Cells.Item(25, 4).Value = "'0.1244"
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value = '  -5.46%  '
$ws.Cells.Item(26, 4).Value = "'7.538"
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 5).Value = '  -4.40%  '
$ws.Cells.Item(27, 4).Value = "'15.69"
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).Value = '  -0.06%  '
$ws.Cells.Item(28, 4).Value = "'0.06298"
$ws.Cells.Item(28, 4).ClearFormats()
$ws.Cells.Item(28, 5).Value = '  -0.70%  '
$ws.Cells.Item(29, 4).Value = "'1.341"
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).Value = '  -4.98%  '
$ws.Cells.Item(30, 5).Value = '  -4.17%  '
$ws.Cells.Item(31, 4).Value = "'3.466"
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).Value = '  -3.01%  '
$ws.Cells.Item(32, 4).Value = "'3.427"
$ws.Cells.Item(32, 4).ClearFormats()
$ws.Cells.Item(32, 5).Value = '  -4.21%  '
$ws.Cells.Item(33, 4).Value = "'1.611"
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).Value = '  -3.73%  '
$ws.Cells.Item(34, 4).Value = "'0.9921"
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = '  -4.59%  '
$ws.Cells.Item(35, 4).Value = "'2.409"
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(35, 5).Value = '  -0.09%  '
$ws.Cells.Item(36, 4).Value = "'0.6016"
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(37, 4).Value = "'2.705"
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).Value = '  -0.49%  '
$ws.Cells.Item(38, 4).Value = "'6.074"
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(38, 5).Value = '  -2.52%  '
$ws.Cells.Item(39, 4).Value = "'0.01608"
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).Value = '  -1.20%  '
$ws.Cells.Item(40, 4).Value = "'1.092.16"
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).Value = '  -2.01%  '
$ws.Cells.Item(41, 4).Value = "'0.8602"
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(42, 5).Value = '  -1.12%  '
$ws.Cells.Item(43, 4).Value = "'99.48"
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).Value = '  -1.69%  '
$ws.Cells.Item(44, 5).Value = '  -1.23%  '
$ws.Cells.Item(45, 4).Value = "'0.00000000110"
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).Value = '  -0.69%  '
$ws.Cells.Item(46, 4).Value = "'56.03"
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).Value = '  -2.13%  '
$ws.Cells.Item(47, 4).Value = "'1.005"
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(47, 5).Value = '  -0.22%  '
$ws.Cells.Item(48, 4).Value = "'0.05248"
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(48, 5).Value = '  -0.42%  '
$ws.Cells.Item(49, 4).Value = "'7.911"
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).Value = '  -3.83%  '
$ws.Cells.Item(50, 4).Value = "'0.4261"
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).Value = '  -1.03%  '
$ws.Cells.Item(51, 4).Value = "'5.881"
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).Value = '  -2.45%  '
